$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.034.08'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '1.829.57'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6234'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.59%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '44.57'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.2910'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.77'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07636'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.64%  '
$ws.Range('D13').Value = '1.830.31'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.956'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6643'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009081'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.80%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.999'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('D19').Value = '28.789.54'
$ws.Range('E19').Value = '  -1.27%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '224.38'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.15%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.181'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.78%  '
$ws.Range('B24').Value = 'BinanceUSD'
$ws.Range('C24').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.53'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.58%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.383'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1355'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.493'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.028'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.046'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('B32').Value = 'Toncoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.200'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05190'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.836'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.57%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.152'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7322'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.16%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.611'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.56%  '
$ws.Range('B38').Value = 'Maker'
$ws.Range('C38').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D38').Value = '1.289.37'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.755'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.84%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01779'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.379'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.29%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8915'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.17%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.50'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.979.69'
$ws.Range('E45').Value = '  +0.27%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5115'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.46%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '63.47'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.58%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000119'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.54%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3971'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.889'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.86%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.643'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.01%  '
